# repull data, push all data, mean calculation
# Updates column F (dSF) values for the rows where the source data changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -5
    8  = 2
    10 = -3
    11 = 8
    16 = 1
    17 = 0
    20 = 2
    27 = 2
    32 = -1
    33 = -2
    37 = 1
    41 = 2
    45 = 2
    48 = 2
    59 = -2
    60 = 1
    64 = 0
    68 = 0
    73 = 0
    75 = -8
    76 = 0
    78 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
